$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Cells.Item(9, 8).Value = 2499
$ws.Cells.Item(9, 9).Value = 1875
$ws.Cells.Item(9, 10).Value = 4995
$ws.Cells.Item(9, 11).Value = 1875
$ws.Cells.Item(9, 12).Value = 4995
$ws.Cells.Item(9, 13).Value = -1706
$ws.Cells.Item(9, 14).Value = -5333

# Row 62
$ws.Cells.Item(62, 8).Value = 5066.3335
$ws.Cells.Item(62, 9).Value = 4699.75
$ws.Cells.Item(62, 11).Value = 4699.75
$ws.Cells.Item(62, 13).Value = -4075.75

# Row 65
$ws.Cells.Item(65, 8).Value = 5066.3335
$ws.Cells.Item(65, 9).Value = 4699.75
$ws.Cells.Item(65, 11).Value = 23498.75
$ws.Cells.Item(65, 13).Value = -20378.75

# Row 70
$ws.Cells.Item(70, 8).Value = 1058.1666
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 1058.1666
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 3174.4998
$ws.Cells.Item(70, 13).ClearContents()
$ws.Cells.Item(70, 14).Value = -3714.4998

# Row 73
$ws.Cells.Item(73, 8).Value = 1058.1666
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 1058.1666
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 3174.4998
$ws.Cells.Item(73, 13).ClearContents()
$ws.Cells.Item(73, 14).Value = -5046.4998

# Row 132
$ws.Cells.Item(132, 8).Value = 1158.8
$ws.Cells.Item(132, 9).Value = 1283.5555
$ws.Cells.Item(132, 10).Value = 36
$ws.Cells.Item(132, 11).Value = 3850.6665
$ws.Cells.Item(132, 12).Value = 108
$ws.Cells.Item(132, 13).Value = -1320.6665
$ws.Cells.Item(132, 14).Value = -5168

# Row 135
$ws.Cells.Item(135, 8).Value = 2916.1
$ws.Cells.Item(135, 9).Value = 3560.8572
$ws.Cells.Item(135, 11).Value = 32047.7148
$ws.Cells.Item(135, 13).Value = -29512.7148

# Row 137
$ws.Cells.Item(137, 8).Value = 3221.5557
$ws.Cells.Item(137, 9).Value = 1799
$ws.Cells.Item(137, 11).Value = 5397
$ws.Cells.Item(137, 13).Value = -2847

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 2978.5833
$ws.Cells.Item(32, 9).Value = 2448.3022
$ws.Cells.Item(32, 10).Value = 7539
$ws.Cells.Item(32, 11).Value = 2448.3022
$ws.Cells.Item(32, 12).Value = 7539
$ws.Cells.Item(32, 13).Value = -2161.3022
$ws.Cells.Item(32, 14).Value = -8113

# Row 61
$ws.Cells.Item(61, 8).Value = 3153.4614
$ws.Cells.Item(61, 9).Value = 2999.818
$ws.Cells.Item(61, 10).Value = 3998.5
$ws.Cells.Item(61, 11).Value = 2999.818
$ws.Cells.Item(61, 12).Value = 3998.5
$ws.Cells.Item(61, 13).Value = -2787.818
$ws.Cells.Item(61, 14).Value = -4422.5

# Row 74
$ws.Cells.Item(74, 8).Value = 1417.7
$ws.Cells.Item(74, 9).Value = 772.125
$ws.Cells.Item(74, 10).Value = 4000
$ws.Cells.Item(74, 11).Value = 772.125
$ws.Cells.Item(74, 12).Value = 4000
$ws.Cells.Item(74, 13).Value = 101.875
$ws.Cells.Item(74, 14).Value = -5748

# Row 77
$ws.Cells.Item(77, 8).Value = 1417.7
$ws.Cells.Item(77, 9).Value = 772.125
$ws.Cells.Item(77, 10).Value = 4000
$ws.Cells.Item(77, 11).Value = 3860.625
$ws.Cells.Item(77, 12).Value = 20000
$ws.Cells.Item(77, 13).Value = 507.375
$ws.Cells.Item(77, 14).Value = -28736

# Row 102
$ws.Cells.Item(102, 8).Value = 48590.6
$ws.Cells.Item(102, 9).Value = 48992
$ws.Cells.Item(102, 10).Value = 47988.5
$ws.Cells.Item(102, 11).Value = 48992
$ws.Cells.Item(102, 12).Value = 47988.5
$ws.Cells.Item(102, 13).Value = -47370
$ws.Cells.Item(102, 14).Value = -51232.5

# Row 132
$ws.Cells.Item(132, 8).Value = 3111.7917
$ws.Cells.Item(132, 9).Value = 2477.8572
$ws.Cells.Item(132, 11).Value = 7433.571599999999
$ws.Cells.Item(132, 13).Value = -4903.571599999999

# Row 136
$ws.Cells.Item(136, 8).Value = 3153.4614
$ws.Cells.Item(136, 9).Value = 2999.818
$ws.Cells.Item(136, 10).Value = 3998.5
$ws.Cells.Item(136, 11).Value = 8999.454000000002
$ws.Cells.Item(136, 12).Value = 11995.5
$ws.Cells.Item(136, 13).Value = -6449.454000000002
$ws.Cells.Item(136, 14).Value = -17095.5

$ws = $wb.Worksheets.Item("BSM")
# Row 76
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).ClearContents()

# Row 79
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).ClearContents()

# Row 134
$ws.Cells.Item(134, 8).Value = 3451
$ws.Cells.Item(134, 9).Value = 3426.0454
$ws.Cells.Item(134, 10).Value = 4000
$ws.Cells.Item(134, 11).Value = 10278.1362
$ws.Cells.Item(134, 12).Value = 12000
$ws.Cells.Item(134, 13).Value = -7743.136200000001
$ws.Cells.Item(134, 14).Value = -17070

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Cells.Item(134, 8).Value = 5739.9287
$ws.Cells.Item(134, 9).Value = 6027.615
$ws.Cells.Item(134, 10).Value = 2000
$ws.Cells.Item(134, 11).Value = 18082.845
$ws.Cells.Item(134, 12).Value = 6000
$ws.Cells.Item(134, 13).Value = -15547.845
$ws.Cells.Item(134, 14).Value = -11070

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Cells.Item(13, 8).Value = 167.85715
$ws.Cells.Item(13, 9).Value = 357.5
$ws.Cells.Item(13, 10).Value = 92
$ws.Cells.Item(13, 11).Value = 1072.5
$ws.Cells.Item(13, 12).Value = 276
$ws.Cells.Item(13, 13).Value = -904.5
$ws.Cells.Item(13, 14).Value = -612

# Row 34
$ws.Cells.Item(34, 8).Value = 1978
$ws.Cells.Item(34, 10).Value = 2795
$ws.Cells.Item(34, 12).Value = 8385
$ws.Cells.Item(34, 14).Value = -8553

# Row 44
$ws.Cells.Item(44, 8).Value = 550.5
$ws.Cells.Item(44, 9).Value = 550.5
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 1651.5
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = -1253.5
$ws.Cells.Item(44, 14).ClearContents()

# Row 60
$ws.Cells.Item(60, 8).Value = 849.1429000000001
$ws.Cells.Item(60, 9).Value = 788.8
$ws.Cells.Item(60, 10).Value = 1000
$ws.Cells.Item(60, 11).Value = 2366.4
$ws.Cells.Item(60, 12).Value = 3000
$ws.Cells.Item(60, 13).Value = -2115.4
$ws.Cells.Item(60, 14).Value = -3502

# Row 113
$ws.Cells.Item(113, 8).Value = 621
$ws.Cells.Item(113, 9).Value = 535.8333
$ws.Cells.Item(113, 10).Value = 723.2
$ws.Cells.Item(113, 11).Value = 1607.4999
$ws.Cells.Item(113, 12).Value = 2169.6
$ws.Cells.Item(113, 13).Value = 562.5001
$ws.Cells.Item(113, 14).Value = -6509.6

# Row 122
$ws.Cells.Item(122, 8).Value = 3746.4783
$ws.Cells.Item(122, 9).Value = 1169.5
$ws.Cells.Item(122, 10).Value = 3863.6135
$ws.Cells.Item(122, 11).Value = 10525.5
$ws.Cells.Item(122, 12).Value = 34772.5215
$ws.Cells.Item(122, 13).Value = -8075.5
$ws.Cells.Item(122, 14).Value = -39672.5215

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Cells.Item(97, 8).Value = 1670
$ws.Cells.Item(97, 9).Value = 1670
$ws.Cells.Item(97, 11).Value = 1670
$ws.Cells.Item(97, 13).Value = -1174

# Row 113
$ws.Cells.Item(113, 8).Value = 1394
$ws.Cells.Item(113, 9).Value = 1395
$ws.Cells.Item(113, 11).Value = 1395
$ws.Cells.Item(113, 13).Value = 775

# Row 132
$ws.Cells.Item(132, 8).Value = 4859.174
$ws.Cells.Item(132, 9).Value = 4790.643
$ws.Cells.Item(132, 10).Value = 4965.778
$ws.Cells.Item(132, 11).Value = 14371.929
$ws.Cells.Item(132, 12).Value = 14897.334
$ws.Cells.Item(132, 13).Value = -11841.929
$ws.Cells.Item(132, 14).Value = -19957.334

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 3249.3333
$ws.Cells.Item(7, 9).Value = 2874.25
$ws.Cells.Item(7, 11).Value = 2874.25
$ws.Cells.Item(7, 13).Value = -2762.25

# Row 40
$ws.Cells.Item(40, 8).Value = 3286.5
$ws.Cells.Item(40, 9).Value = 3398.8333
$ws.Cells.Item(40, 11).Value = 3398.8333
$ws.Cells.Item(40, 13).Value = -3262.8333

# Row 61
$ws.Cells.Item(61, 8).Value = 4992.25
$ws.Cells.Item(61, 9).Value = 4990
$ws.Cells.Item(61, 11).Value = 4990
$ws.Cells.Item(61, 13).Value = -4788

# Row 113
$ws.Cells.Item(113, 8).Value = 4992.25
$ws.Cells.Item(113, 9).Value = 4990
$ws.Cells.Item(113, 11).Value = 4990
$ws.Cells.Item(113, 13).Value = -2820

# Row 126
$ws.Cells.Item(126, 8).Value = 3249.3333
$ws.Cells.Item(126, 9).Value = 2874.25
$ws.Cells.Item(126, 11).Value = 8622.75
$ws.Cells.Item(126, 13).Value = -6152.75

# Row 132
$ws.Cells.Item(132, 8).Value = 5747.125
$ws.Cells.Item(132, 9).Value = 5002
$ws.Cells.Item(132, 11).Value = 15006
$ws.Cells.Item(132, 13).Value = -12476

# Row 136
$ws.Cells.Item(136, 8).Value = 23464.87
$ws.Cells.Item(136, 9).Value = 2770.2
$ws.Cells.Item(136, 11).Value = 8310.599999999999
$ws.Cells.Item(136, 13).Value = -5760.599999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Cells.Item(14, 8).Value = 725
$ws.Cells.Item(14, 10).Value = 1500
$ws.Cells.Item(14, 12).Value = 1500
$ws.Cells.Item(14, 14).Value = -1836

# Row 113
$ws.Cells.Item(113, 8).Value = 759.4286
$ws.Cells.Item(113, 9).Value = 795.4
$ws.Cells.Item(113, 11).Value = 2386.2
$ws.Cells.Item(113, 13).Value = -216.1999999999998

# Row 126
$ws.Cells.Item(126, 8).Value = 2027.1428
$ws.Cells.Item(126, 9).Value = 1838.2
$ws.Cells.Item(126, 11).Value = 5514.6
$ws.Cells.Item(126, 13).Value = -3044.6
